$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.654227
$ws.Cells.Item(2, 8).Value = 4.962681
$ws.Cells.Item(2, 9).Value = 0.4107585939979205
$ws.Cells.Item(2, 10).Value = 0.4107585939979205
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.168007333333333
$ws.Cells.Item(2, 14).Value = 3.504022
$ws.Cells.Item(2, 15).Value = 0.1638609704511517
$ws.Cells.Item(2, 16).Value = 0.1638609704511517
$ws.Cells.Item(2, 17).Value = 1.932149266998
$ws.Cells.Item(2, 18).Value = 17.389343402982
$ws.Cells.Item(2, 19).Value = 0.06730730183364989
$ws.Cells.Item(2, 20).Value = 0.06730730183364989

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.654227
$ws.Cells.Item(3, 8).Value = 4.962681
$ws.Cells.Item(3, 9).Value = 0.4107585939979205
$ws.Cells.Item(3, 10).Value = 0.4107585939979205
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.578098999999999
$ws.Cells.Item(3, 14).Value = 10.734297
$ws.Cells.Item(3, 15).Value = 0.5019752511630595
$ws.Cells.Item(3, 16).Value = 0.5019752511630595
$ws.Cells.Item(3, 17).Value = 5.918987974472999
$ws.Cells.Item(3, 18).Value = 53.27089177025699
$ws.Cells.Item(3, 19).Value = 0.2061906483894914
$ws.Cells.Item(3, 20).Value = 0.2061906483894914

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.654227
$ws.Cells.Item(4, 8).Value = 4.962681
$ws.Cells.Item(4, 9).Value = 0.4107585939979205
$ws.Cells.Item(4, 10).Value = 0.4107585939979205
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.915392333333333
$ws.Cells.Item(4, 14).Value = 5.746177
$ws.Cells.Item(4, 15).Value = 0.2687123938160456
$ws.Cells.Item(4, 16).Value = 0.2687123938160456
$ws.Cells.Item(4, 17).Value = 3.168493713393
$ws.Cells.Item(4, 18).Value = 28.516443420537
$ws.Cells.Item(4, 19).Value = 0.1103759250736944
$ws.Cells.Item(4, 20).Value = 0.1103759250736944

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.654227
$ws.Cells.Item(5, 8).Value = 4.962681
$ws.Cells.Item(5, 9).Value = 0.4107585939979205
$ws.Cells.Item(5, 10).Value = 0.4107585939979205
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.46654
$ws.Cells.Item(5, 14).Value = 1.39962
$ws.Cells.Item(5, 15).Value = 0.06545138456974327
$ws.Cells.Item(5, 16).Value = 0.06545138456974327
$ws.Cells.Item(5, 17).Value = 0.7717630645799999
$ws.Cells.Item(5, 18).Value = 6.945867581219999
$ws.Cells.Item(5, 19).Value = 0.02688471870108494
$ws.Cells.Item(5, 20).Value = 0.02688471870108494

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.288726
$ws.Cells.Item(6, 8).Value = 3.866178
$ws.Cells.Item(6, 9).Value = 0.3200015957958394
$ws.Cells.Item(6, 10).Value = 0.3200015957958394
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.168007333333333
$ws.Cells.Item(6, 14).Value = 3.504022
$ws.Cells.Item(6, 15).Value = 0.1638609704511517
$ws.Cells.Item(6, 16).Value = 0.1638609704511517
$ws.Cells.Item(6, 17).Value = 1.505241418657333
$ws.Cells.Item(6, 18).Value = 13.547172767916
$ws.Cells.Item(6, 19).Value = 0.05243577203302344
$ws.Cells.Item(6, 20).Value = 0.05243577203302344

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.288726
$ws.Cells.Item(7, 8).Value = 3.866178
$ws.Cells.Item(7, 9).Value = 0.3200015957958394
$ws.Cells.Item(7, 10).Value = 0.3200015957958394
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.578098999999999
$ws.Cells.Item(7, 14).Value = 10.734297
$ws.Cells.Item(7, 15).Value = 0.5019752511630595
$ws.Cells.Item(7, 16).Value = 0.5019752511630595
$ws.Cells.Item(7, 17).Value = 4.611189211873999
$ws.Cells.Item(7, 18).Value = 41.50070290686599
$ws.Cells.Item(7, 19).Value = 0.1606328814221963
$ws.Cells.Item(7, 20).Value = 0.1606328814221963

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.288726
$ws.Cells.Item(8, 8).Value = 3.866178
$ws.Cells.Item(8, 9).Value = 0.3200015957958394
$ws.Cells.Item(8, 10).Value = 0.3200015957958394
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.915392333333333
$ws.Cells.Item(8, 14).Value = 5.746177
$ws.Cells.Item(8, 15).Value = 0.2687123938160456
$ws.Cells.Item(8, 16).Value = 0.2687123938160456
$ws.Cells.Item(8, 17).Value = 2.468415900167333
$ws.Cells.Item(8, 18).Value = 22.215743101506
$ws.Cells.Item(8, 19).Value = 0.08598839483125463
$ws.Cells.Item(8, 20).Value = 0.08598839483125464

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.288726
$ws.Cells.Item(9, 8).Value = 3.866178
$ws.Cells.Item(9, 9).Value = 0.3200015957958394
$ws.Cells.Item(9, 10).Value = 0.3200015957958394
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.46654
$ws.Cells.Item(9, 14).Value = 1.39962
$ws.Cells.Item(9, 15).Value = 0.06545138456974327
$ws.Cells.Item(9, 16).Value = 0.06545138456974327
$ws.Cells.Item(9, 17).Value = 0.6012422280399998
$ws.Cells.Item(9, 18).Value = 5.411180052359999
$ws.Cells.Item(9, 19).Value = 0.02094454750936503
$ws.Cells.Item(9, 20).Value = 0.02094454750936503

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8858993333333333
$ws.Cells.Item(10, 8).Value = 2.657698
$ws.Cells.Item(10, 9).Value = 0.2199763179924491
$ws.Cells.Item(10, 10).Value = 0.2199763179924491
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.168007333333333
$ws.Cells.Item(10, 14).Value = 3.504022
$ws.Cells.Item(10, 15).Value = 0.1638609704511517
$ws.Cells.Item(10, 16).Value = 0.1638609704511517
$ws.Cells.Item(10, 17).Value = 1.034736917928444
$ws.Cells.Item(10, 18).Value = 9.312632261355999
$ws.Cells.Item(10, 19).Value = 0.03604553294251386
$ws.Cells.Item(10, 20).Value = 0.03604553294251386

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.8858993333333333
$ws.Cells.Item(11, 8).Value = 2.657698
$ws.Cells.Item(11, 9).Value = 0.2199763179924491
$ws.Cells.Item(11, 10).Value = 0.2199763179924491
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 3.578098999999999
$ws.Cells.Item(11, 14).Value = 10.734297
$ws.Cells.Item(11, 15).Value = 0.5019752511630595
$ws.Cells.Item(11, 16).Value = 0.5019752511630595
$ws.Cells.Item(11, 17).Value = 3.169835518700666
$ws.Cells.Item(11, 18).Value = 28.52851966830599
$ws.Cells.Item(11, 19).Value = 0.1104226674741847
$ws.Cells.Item(11, 20).Value = 0.1104226674741847

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.8858993333333333
$ws.Cells.Item(12, 8).Value = 2.657698
$ws.Cells.Item(12, 9).Value = 0.2199763179924491
$ws.Cells.Item(12, 10).Value = 0.2199763179924491
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.915392333333333
$ws.Cells.Item(12, 14).Value = 5.746177
$ws.Cells.Item(12, 15).Value = 0.2687123938160456
$ws.Cells.Item(12, 16).Value = 0.2687123938160456
$ws.Cells.Item(12, 17).Value = 1.696844791171778
$ws.Cells.Item(12, 18).Value = 15.271603120546
$ws.Cells.Item(12, 19).Value = 0.05911036299059066
$ws.Cells.Item(12, 20).Value = 0.05911036299059067

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.8858993333333333
$ws.Cells.Item(13, 8).Value = 2.657698
$ws.Cells.Item(13, 9).Value = 0.2199763179924491
$ws.Cells.Item(13, 10).Value = 0.2199763179924491
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.46654
$ws.Cells.Item(13, 14).Value = 1.39962
$ws.Cells.Item(13, 15).Value = 0.06545138456974327
$ws.Cells.Item(13, 16).Value = 0.06545138456974327
$ws.Cells.Item(13, 17).Value = 0.4133074749733333
$ws.Cells.Item(13, 18).Value = 3.71976727476
$ws.Cells.Item(13, 19).Value = 0.01439775458515992
$ws.Cells.Item(13, 20).Value = 0.01439775458515992

# Row 14
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.1983963333333333
$ws.Cells.Item(14, 8).Value = 0.595189
$ws.Cells.Item(14, 9).Value = 0.04926349221379096
$ws.Cells.Item(14, 10).Value = 0.04926349221379096
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.168007333333333
$ws.Cells.Item(14, 14).Value = 3.504022
$ws.Cells.Item(14, 15).Value = 0.1638609704511517
$ws.Cells.Item(14, 16).Value = 0.1638609704511517
$ws.Cells.Item(14, 17).Value = 0.2317283722397778
$ws.Cells.Item(14, 18).Value = 2.085555350158
$ws.Cells.Item(14, 19).Value = 0.008072363641964543
$ws.Cells.Item(14, 20).Value = 0.008072363641964543

# Row 15
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.1983963333333333
$ws.Cells.Item(15, 8).Value = 0.595189
$ws.Cells.Item(15, 9).Value = 0.04926349221379096
$ws.Cells.Item(15, 10).Value = 0.04926349221379096
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 3.578098999999999
$ws.Cells.Item(15, 14).Value = 10.734297
$ws.Cells.Item(15, 15).Value = 0.5019752511630595
$ws.Cells.Item(15, 16).Value = 0.5019752511630595
$ws.Cells.Item(15, 17).Value = 0.7098817219036665
$ws.Cells.Item(15, 18).Value = 6.388935497132999
$ws.Cells.Item(15, 19).Value = 0.02472905387718714
$ws.Cells.Item(15, 20).Value = 0.02472905387718714

# Row 16
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.1983963333333333
$ws.Cells.Item(16, 8).Value = 0.595189
$ws.Cells.Item(16, 9).Value = 0.04926349221379096
$ws.Cells.Item(16, 10).Value = 0.04926349221379096
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 1.915392333333333
$ws.Cells.Item(16, 14).Value = 5.746177
$ws.Cells.Item(16, 15).Value = 0.2687123938160456
$ws.Cells.Item(16, 16).Value = 0.2687123938160456
$ws.Cells.Item(16, 17).Value = 0.3800068158281111
$ws.Cells.Item(16, 18).Value = 3.420061342453
$ws.Cells.Item(16, 19).Value = 0.01323771092050589
$ws.Cells.Item(16, 20).Value = 0.01323771092050589

# Row 17
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.1983963333333333
$ws.Cells.Item(17, 8).Value = 0.595189
$ws.Cells.Item(17, 9).Value = 0.04926349221379096
$ws.Cells.Item(17, 10).Value = 0.04926349221379096
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.46654
$ws.Cells.Item(17, 14).Value = 1.39962
$ws.Cells.Item(17, 15).Value = 0.06545138456974327
$ws.Cells.Item(17, 16).Value = 0.06545138456974327
$ws.Cells.Item(17, 17).Value = 0.09255982535333332
$ws.Cells.Item(17, 18).Value = 0.83303842818
$ws.Cells.Item(17, 19).Value = 0.003224363774133385
$ws.Cells.Item(17, 20).Value = 0.003224363774133385
